$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New detector-hall parameter rows appended below the existing table
$ws.Range("A7").Value = "Distance from end of straighg to detector hall"
$ws.Range("B7").Value = "HallWallDist"
$ws.Range("C7").Value = 50
$ws.Range("D7").Value = "m"
$ws.Range("E7").Value = "nuSIM-2021-01"

$ws.Range("B8").Value = "DetHlfWdth"
$ws.Range("A8").Value = "Detector half width"
$ws.Range("C8").Value = 2
$ws.Range("D8").Value = "m"
$ws.Range("E8").Value = "nuSIM-2021-01"

$ws.Range("A9").Value = "Detector length"
$ws.Range("B9").Value = "DetLngth"
$ws.Range("C9").Value = 4
$ws.Range("D9").Value = "m"
$ws.Range("E9").Value = "nuSIM-2021-01"

$ws.Range("A10").Value = "Distance from hall wall to detector entrance"
$ws.Range("B10").Value = "Hall2Det"
$ws.Range("C10").Value = 5
$ws.Range("D10").Value = "m"
$ws.Range("E10").Value = "nuSIM-2021-01"

# Column A widens to fit the new, longer parameter descriptions (closest
# reachable value to the target 39.33203125 given this runtime's width
# quantization)
$ws.Columns.Item(1).ColumnWidth = 38.5

# Selection moves, matching the author's final cursor position
$ws.Range("N22").Select()
